$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 650
$ws.Range("F3").Value = 6175
$ws.Range("F7").Value = 1022
$ws.Range("F8").Value = 411
$ws.Range("F9").Value = 1390
$ws.Range("F10").Value = 3145
$ws.Range("F11").Value = 458
$ws.Range("F12").Value = 1999
$ws.Range("F15").Value = 209
$ws.Range("F16").Value = 93
$ws.Range("F17").Value = 194
$ws.Range("F18").Value = 1013
$ws.Range("F21").Value = 117
$ws.Range("F22").Value = 3778
$ws.Range("F23").Value = 1197
$ws.Range("F24").Value = 2997
$ws.Range("F26").Value = 2549
$ws.Range("F27").Value = 4347
$ws.Range("F28").Value = 114
$ws.Range("F29").Value = 940
$ws.Range("F31").Value = 2258
$ws.Range("F32").Value = 151
$ws.Range("F34").Value = 57
$ws.Range("F36").Value = 38
$ws.Range("F37").Value = 1052
$ws.Range("F38").Value = 1303
$ws.Range("F40").Value = 1143
$ws.Range("F41").Value = 737
$ws.Range("F42").Value = 658
$ws.Range("F43").Value = 446
$ws.Range("F44").Value = 28
$ws.Range("F45").Value = 136
$ws.Range("F47").Value = 13
$ws.Range("F48").Value = 328
$ws.Range("F49").Value = 3631

# Sheet: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F14").Value = 3
$ws.Range("F24").Value = 35

# Sheet: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 6175
$ws.Range("F6").Value = 411
$ws.Range("F7").Value = 1390
$ws.Range("F8").Value = 3145
$ws.Range("F9").Value = 458
$ws.Range("F11").Value = 1999
$ws.Range("F15").Value = 209
$ws.Range("F18").Value = 93
$ws.Range("F19").Value = 194
$ws.Range("F20").Value = 1013
$ws.Range("F21").Value = 3
$ws.Range("F23").Value = 117
$ws.Range("F24").Value = 3778
$ws.Range("F26").Value = 1197
$ws.Range("F28").Value = 2997
$ws.Range("F29").Value = 2549
$ws.Range("F30").Value = 4347
$ws.Range("F31").Value = 114
$ws.Range("F32").Value = 940
$ws.Range("F33").Value = 2261
$ws.Range("F34").Value = 38
$ws.Range("F35").Value = 1052
$ws.Range("F36").Value = 1303
$ws.Range("F38").Value = 1143
$ws.Range("F39").Value = 737
$ws.Range("F41").Value = 446
$ws.Range("F43").Value = 35
$ws.Range("F44").Value = 28
$ws.Range("F46").Value = 136
$ws.Range("F47").Value = 328
$ws.Range("F48").Value = 3631
